# Update the "Förändrad" (changed) date column C for rows 2-12
# from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = 45224
}
